$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.710.74"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "3.503.51"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "585.91"
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("D6").Value = "132.27"
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("D7").Value = "3.504.45"
$ws.Range("E7").Value = "  -1.25%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").Value = "4.100.80"
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("D14").Value = "27.70"
$ws.Range("E14").Value = "  +3.08%  "
$ws.Range("E15").Value = "  -1.83%  "
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("D17").Value = "3.506.98"
$ws.Range("E17").Value = "  -1.18%  "
$ws.Range("D18").Value = "64.749.65"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").Value = "9.99"
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("D22").Value = "391.62"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "3.644.51"
$ws.Range("E24").Value = "  -1.37%  "
$ws.Range("D25").Value = "74.10"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -4.72%  "
$ws.Range("D28").Value = "1.58"
$ws.Range("E28").Value = "  +1.75%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "7.40"
$ws.Range("E30").Value = "  -4.80%  "
$ws.Range("D31").Value = "2.25"
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("D32").Value = "8.19"
$ws.Range("E32").Value = "  -3.94%  "
$ws.Range("D33").Value = "3.505.79"
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("E36").Value = "  -0.67%  "
$ws.Range("E37").Value = "  +3.74%  "
$ws.Range("D38").Value = "171.99"
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("D39").Value = "5.19"
$ws.Range("E39").Value = "  +3.49%  "
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").Value = "0.0809"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  -1.46%  "
$ws.Range("D43").Value = "26.25"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").Value = "42.34"
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("E46").Value = "  -2.30%  "
$ws.Range("D47").Value = "4.39"
$ws.Range("E47").Value = "  -0.95%  "
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("D49").Value = "2.479.68"
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("D51").Value = "0.903"
$ws.Range("E51").Value = "  +3.77%  "
